$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: remove old rows 2-5 (Potion/Scroll/Gold/Ammo item-type charge rows).
# This shifts Food(old r6)->r2, Bandage(old r7)->r3, Artifact(old r8)->r4,
# Dungeonpedia(old r9)->r5, Blueprint(old r10)->r6.
$ws.Rows("2:5").Delete()

# Step 2: insert two fresh blank rows at position 2 so the Artifact row (now at row 4)
# shifts back down to row 6 (preserving its autofit row height of 171), giving:
# r2=blank, r3=blank, r4=Food(old r6), r5=Bandage(old r7), r6=Artifact(old r8),
# r7=Dungeonpedia(old r9), r8=Blueprint(old r10)
$ws.Rows("2:3").Insert()

# Step 3: populate the four new "loot card" rows (2-5) with the new item-type data
$ws.Range("A2").Value = "材料包"
$ws.Range("B2").Value = "为1张道具牌充1能。可以额外消耗1张同名牌，改为充3能。"
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = "Food"
$ws.Range("E2").Value = "① Activate: Add 1 food. Add 1 extra potion for every 2 [Food] card under this card."
$ws.Range("D2").Style = "常规 2"
$ws.Range("E2").Style = "常规 2"

$ws.Range("A3").Value = "绷带"
$ws.Range("B3").Value = "回复1生命。可以额外消耗1张同名牌，改为回复3生命。"
$ws.Range("C3").Value = 4
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("D3").Style = "常规 2"
$ws.Range("E3").Style = "常规 2"

$ws.Range("A4").Value = "壶"
$ws.Range("B4").Value = "从主牌堆抽2张牌。每次整理战利品只能使用1张。"
$ws.Range("C4").Value = 4
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("D4").Style = "常规 2"
$ws.Range("E4").Style = "常规 2"

$ws.Range("A5").Value = "魔镜"
$ws.Range("B5").Value = "选手牌或场上1张牌，从主牌堆获得其1张同名牌。"
$ws.Range("C5").Value = 4
$ws.Range("D5").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("D5").Style = "常规 2"
$ws.Range("E5").Style = "常规 2"

# Step 4: set the active cell/selection like the saved workbook
$ws.Range("F6").Select() | Out-Null

